$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44424,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',30,32000,32000,32000,'$/caja 18 kilos','Perú',1778,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44305,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',25,30000,30000,30000,'$/caja 18 kilos','Perú',1667,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44326,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Especial',16,35000,35000,35000,'$/caja 18 kilos','Perú',1944,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44326,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',25,30000,30000,30000,'$/caja 18 kilos','Perú',1667,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44326,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Segunda',20,28000,28000,28000,'$/caja 18 kilos','Perú',1556,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44270,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Especial',70,38000,38000,38000,'$/caja 18 kilos','Perú',2111,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44389,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',35,29000,29000,29000,'$/caja 18 kilos','Perú',1611,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44389,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Segunda',20,27000,27000,27000,'$/caja 18 kilos','Perú',1500,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44382,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Especial',20,35000,35000,35000,'$/caja 18 kilos','Perú',1944,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44382,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',30,32000,32000,32000,'$/caja 18 kilos','Perú',1778,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44382,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Segunda',15,30000,30000,30000,'$/caja 18 kilos','Perú',1667,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44445,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',30,32000,32000,32000,'$/caja 18 kilos','Perú',1778,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44403,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Especial',25,33000,33000,33000,'$/caja 18 kilos','Perú',1833,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44403,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',45,30000,30000,30000,'$/caja 18 kilos','Perú',1667,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44403,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Segunda',15,28000,28000,28000,'$/caja 18 kilos','Perú',1556,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44354,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',40,38000,38000,38000,'$/caja 18 kilos','Perú',2111,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44333,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',30,38000,38000,38000,'$/caja 18 kilos','Perú',2111,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44333,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Segunda',25,35000,35000,35000,'$/caja 18 kilos','Perú',1944,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44431,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',30,32000,32000,32000,'$/caja 18 kilos','Perú',1778,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44312,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',160,26000,26000,26000,'$/caja 18 kilos','Perú',1444,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44410,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Especial',15,32000,32000,32000,'$/caja 18 kilos','Perú',1778,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44410,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',25,30000,30000,30000,'$/caja 18 kilos','Perú',1667,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44410,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Segunda',10,28000,28000,28000,'$/caja 18 kilos','Perú',1556,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44284,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',120,23000,23000,23000,'$/caja 18 kilos','Perú',1278,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44396,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',35,37000,37000,37000,'$/caja 18 kilos','Perú',2056,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44396,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Segunda',15,34000,34000,34000,'$/caja 18 kilos','Perú',1889,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44340,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',35,37000,37000,37000,'$/caja 18 kilos','Perú',2056,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44340,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Segunda',20,35000,35000,35000,'$/caja 18 kilos','Perú',1944,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44277,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',100,30000,30000,30000,'$/caja 18 kilos','Perú',1667,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44277,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Segunda',60,28000,28000,28000,'$/caja 18 kilos','Perú',1556,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44291,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',30,25000,25000,25000,'$/caja 18 kilos','Perú',1389,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44417,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',15,28000,28000,28000,'$/caja 18 kilos','Perú',1556,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44438,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',30,32000,32000,32000,'$/caja 18 kilos','Perú',1778,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44435,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',30,32000,32000,32000,'$/caja 18 kilos','Perú',1778,18),
    @(9,'Vega Central Mapocho de Santiago','Metropolitana',44319,13,'Fruta',100108,'Tropicales y subtropicales',100108003,'Maracuyá','Sin especificar','Primera',140,27000,27000,27000,'$/caja 18 kilos','Perú',1500,18)
)

$dateNumberFormat = $ws.Cells.Item(2, 4).NumberFormat

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $colIndex = $j + 1
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowValues[$j]
    }
    $ws.Cells.Item($rowIndex, 4).NumberFormat = $dateNumberFormat
}
